# Auto-generated script to apply Maduin_Profits price-data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value2 = 238.35715
$ws.Range("I5").Value2 = 116.77778
$ws.Range("J5").Value2 = 457.2
$ws.Range("K5").Value2 = 116.77778
$ws.Range("L5").Value2 = 457.2
$ws.Range("M5").Value2 = -1.777780000000007
$ws.Range("N5").Value2 = -687.2
$ws.Range("H10").Value2 = 100
$ws.Range("I10").Value2 = 100
$ws.Range("K10").Value2 = 100
$ws.Range("M10").Value2 = 193
$ws.Range("H15").Value2 = 575.53845
$ws.Range("I15").Value2 = 575.53845
$ws.Range("K15").Value2 = 1726.61535
$ws.Range("M15").Value2 = -1557.61535
$ws.Range("H21").Value2 = 1750
$ws.Range("I21").Value2 = 1750
$ws.Range("K21").Value2 = 1750
$ws.Range("M21").Value2 = -1282
$ws.Range("H23").Value2 = 1750
$ws.Range("I23").Value2 = 1750
$ws.Range("K23").Value2 = 1750
$ws.Range("M23").Value2 = -1516
$ws.Range("H33").Value2 = 337
$ws.Range("I33").Value2 = 276
$ws.Range("K33").Value2 = 276
$ws.Range("M33").Value2 = -47
$ws.Range("H62").Value2 = 7149.143
$ws.Range("I62").Value2 = 6220.4
$ws.Range("K62").Value2 = 6220.4
$ws.Range("M62").Value2 = -5596.4
$ws.Range("H63").Value2 = 0
$ws.Range("J63").Value2 = 0
$ws.Range("L63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value2 = 10363.637
$ws.Range("J64").Value2 = 12000
$ws.Range("L64").Value2 = 12000
$ws.Range("N64").Value2 = -12496
$ws.Range("H65").Value2 = 7149.143
$ws.Range("I65").Value2 = 6220.4
$ws.Range("K65").Value2 = 31102
$ws.Range("M65").Value2 = -27982
$ws.Range("H66").Value2 = 0
$ws.Range("J66").Value2 = 0
$ws.Range("L66").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value2 = 10363.637
$ws.Range("J67").Value2 = 12000
$ws.Range("L67").Value2 = 12000
$ws.Range("N67").Value2 = -13716
$ws.Range("H113").Value2 = 2779
$ws.Range("I113").Value2 = 2400
$ws.Range("J113").Value2 = 3347.5
$ws.Range("K113").Value2 = 2400
$ws.Range("L113").Value2 = 3347.5
$ws.Range("M113").Value2 = 854
$ws.Range("N113").Value2 = -9855.5
$ws.Range("H133").Value2 = 49999
$ws.Range("J133").Value2 = 49999
$ws.Range("L133").Value2 = 49999
$ws.Range("N133").Value2 = -60119

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value2 = 25000
$ws.Range("J43").Value2 = 25000
$ws.Range("L43").Value2 = 25000
$ws.Range("N43").Value2 = -25626
$ws.Range("H122").Value2 = 2733.3333
$ws.Range("I122").Value2 = 4000
$ws.Range("K122").Value2 = 12000
$ws.Range("M122").Value2 = -9550
$ws.Range("H132").Value2 = 1602.4667
$ws.Range("I132").Value2 = 1384.1818
$ws.Range("K132").Value2 = 4152.5454
$ws.Range("M132").Value2 = -1622.5454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value2 = 4520.5557
$ws.Range("I105").Value2 = 4460.625
$ws.Range("K105").Value2 = 4460.625
$ws.Range("M105").Value2 = -2713.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value2 = 2735
$ws.Range("I23").Value2 = 2735
$ws.Range("K23").Value2 = 2735
$ws.Range("M23").Value2 = -2495
$ws.Range("H27").Value2 = 2735
$ws.Range("I27").Value2 = 2735
$ws.Range("K27").Value2 = 2735
$ws.Range("M27").Value2 = -2543
$ws.Range("H99").Value2 = 5755.9375
$ws.Range("I99").Value2 = 5128.846
$ws.Range("J99").Value2 = 8473.333000000001
$ws.Range("K99").Value2 = 5128.846
$ws.Range("L99").Value2 = 8473.333000000001
$ws.Range("M99").Value2 = -3630.846
$ws.Range("N99").Value2 = -11469.333
$ws.Range("H126").Value2 = 5755.9375
$ws.Range("I126").Value2 = 5128.846
$ws.Range("J126").Value2 = 8473.333000000001
$ws.Range("K126").Value2 = 15386.538
$ws.Range("L126").Value2 = 25419.999
$ws.Range("M126").Value2 = -12916.538
$ws.Range("N126").Value2 = -30359.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value2 = 7929.1113
$ws.Range("J80").Value2 = 10643.667
$ws.Range("L80").Value2 = 31931.001
$ws.Range("N80").Value2 = -33803.001
$ws.Range("H83").Value2 = 7929.1113
$ws.Range("J83").Value2 = 10643.667
$ws.Range("L83").Value2 = 95793.003
$ws.Range("N83").Value2 = -105153.003
$ws.Range("H104").Value2 = 5666.3335
$ws.Range("J104").Value2 = 5666.3335
$ws.Range("L104").Value2 = 16999.0005
$ws.Range("N104").Value2 = -22241.0005
$ws.Range("H122").Value2 = 1088.5
$ws.Range("I122").Value2 = 200
$ws.Range("J122").Value2 = 1977
$ws.Range("K122").Value2 = 1800
$ws.Range("L122").Value2 = 17793
$ws.Range("M122").Value2 = 650
$ws.Range("N122").Value2 = -22693
$ws.Range("H131").Value2 = 1037.25
$ws.Range("J131").Value2 = 1066.3334
$ws.Range("L131").Value2 = 3199.0002
$ws.Range("N131").Value2 = -13279.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value2 = 1266.6666
$ws.Range("I6").Value2 = 1266.6666
$ws.Range("K6").Value2 = 1266.6666
$ws.Range("M6").Value2 = -1153.6666
$ws.Range("H16").Value2 = 1266.6666
$ws.Range("I16").Value2 = 1266.6666
$ws.Range("K16").Value2 = 1266.6666
$ws.Range("M16").Value2 = -1016.6666
$ws.Range("H113").Value2 = 905
$ws.Range("I113").Value2 = 1000
$ws.Range("J113").Value2 = 810
$ws.Range("K113").Value2 = 1000
$ws.Range("L113").Value2 = 810
$ws.Range("M113").Value2 = 1170
$ws.Range("N113").Value2 = -5150

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value2 = 15399.5
$ws.Range("I29").Value2 = 800
$ws.Range("K29").Value2 = 800
$ws.Range("M29").Value2 = -505
$ws.Range("H74").Value2 = 69999
$ws.Range("J74").Value2 = 69999
$ws.Range("L74").Value2 = 69999
$ws.Range("N74").Value2 = -71995
$ws.Range("H77").Value2 = 69999
$ws.Range("J77").Value2 = 69999
$ws.Range("L77").Value2 = 209997
$ws.Range("N77").Value2 = -219981
$ws.Range("H132").Value2 = 4999.6665
$ws.Range("I132").Value2 = 5999.5
$ws.Range("K132").Value2 = 17998.5
$ws.Range("M132").Value2 = -15468.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value2 = 45000000
$ws.Range("I3").Value2 = 0
$ws.Range("J3").Value2 = 45000000
$ws.Range("K3").Value2 = 0
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value2 = -45000228
$ws.Range("H30").Value2 = 15000
$ws.Range("J30").Value2 = 15000
$ws.Range("L30").Value2 = 15000
$ws.Range("N30").Value2 = -15214
$ws.Range("H127").Value2 = 0
$ws.Range("I127").Value2 = 0
$ws.Range("K127").Value2 = 0
$ws.Range("M127").ClearContents()
